$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.221991658210754
$ws.Range("B1").Value = 2.693856716156006
$ws.Range("C1").Value = 4.310568809509277
$ws.Range("D1").Value = 2.141778707504272
$ws.Range("E1").Value = 1.159821510314941
